$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-09-03T11:48:25+00:00"

# --- Elements sheet: clear stale Condition(s) refs, fix RIM mapping casing ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI4").Value = ""
$elements.Range("AK4").Value = "n/a"
$elements.Range("AI5").Value = ""
$elements.Range("AI7").Value = ""
$elements.Range("AI9").Value = ""
$elements.Range("AI10").Value = ""
$elements.Range("AI12").Value = ""
$elements.Range("AI14").Value = ""
$elements.Range("AI16").Value = ""
